$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "iron man"
$ws.Range("A4").Value = "hulk"

$ws.Range("A4").Select()
